$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the static "About" message in C13 with a formula that builds the
# same text, padded with 50 trailing spaces (widget text solution).
$ws.Range("C13").Formula = '="Memolab 2021 \nversion a \nAuteur : \nJean-Claude Vouillamoz  \nmail : \njcvouillamoz@gmail.com  \nPhone : \n+ 41 79 212 84 52 \nLicense GNU \n"&REPT(" ",50)'

# Move the active selection to C14, matching the saved cursor position.
$ws.Range("C14").Select()
